$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# First body paragraph: "Kehidupan manusia ... seseorang."
$para1 = $tr.Paragraphs(1)

# Justify the paragraph (pPr algn="just")
$para1.ParagraphFormat.Alignment = 4   # ppAlignJustify

# Prepend a tab character to the start of the paragraph's text
# (only touches the first run, "Kehidupan ")
$lead = $tr.Characters(1, 10)
$lead.Text = "`tKehidupan "

# Second (trailing empty) paragraph: also justify it
$para2 = $tr.Paragraphs(2)
$para2.ParagraphFormat.Alignment = 4   # ppAlignJustify
